$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsBDP = $wb.Worksheets.Item("BDPbES")

# ---------------------------------------------------------------------------
# 1) BDPbES sheet: three new fuel rows (15-17), each priority = 2, with a
#    formula in C:AK that copies across from column B (matching the existing
#    pattern used by every other row on this sheet).
# ---------------------------------------------------------------------------
$wsBDP.Range("A15").Value = "crude oil"
$wsBDP.Range("A16").Value = "heavy or residual fuel oil"
$wsBDP.Range("A17").Value = "municipal solid waste"

$wsBDP.Range("B15").Value = 2
$wsBDP.Range("B16").Value = 2
$wsBDP.Range("B17").Value = 2

$wsBDP.Range("C15:AK17").Formula = "=`$B15"

# ---------------------------------------------------------------------------
# 2) BDPbES sheet: petroleum (row 11) priority changes from 1 to 2 across the
#    whole row (literal values, not formulas).
# ---------------------------------------------------------------------------
$wsBDP.Range("B11:AK11").Value = 2

# ---------------------------------------------------------------------------
# 3) BDPbES sheet: new header label in A1 describing column B onward, bold +
#    wrapped, with a taller header row, and a wider column A to fit it.
# ---------------------------------------------------------------------------
$wsBDP.Range("A1").Value = "Priority Order (dimensionless)"
$wsBDP.Range("A1").Font.Bold = $true
$wsBDP.Range("A1").WrapText = $true
$wsBDP.Rows.Item(1).RowHeight = 30
$wsBDP.Columns.Item(1).ColumnWidth = 27.166666666666668

# ---------------------------------------------------------------------------
# 4) About sheet: new explanatory note (red text) clarifying the HK-specific
#    change to the petroleum priority.
# ---------------------------------------------------------------------------
$wsAbout.Range("A12").Value = "US model petrolium is 1, changed 2 for HK"
$wsAbout.Range("A12").Font.Color = 255

# ---------------------------------------------------------------------------
# 5) Selections / active sheet, matching the saved UI state in the workbook.
# ---------------------------------------------------------------------------
$wsAbout.Range("D18").Select() | Out-Null
$wsBDP.Activate() | Out-Null
$wsBDP.Range("B9").Select() | Out-Null
